$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear existing data rows (keep header row 1)
$ws.Range("A2:G146").ClearContents() | Out-Null

$data = New-Object 'object[,]' 149,7
$data[0,0] = 17
$data[0,1] = 'Dobrrzyńskiej'
$data[0,2] = 'ziemi Dobrrzyńskiej ('
$data[0,3] = 'dobrzyńska'
$data[0,4] = 'Dobrrzyńskiej'
$data[0,5] = 'dobrrzyńska'
$data[0,6] = 'spelling'
$data[1,0] = 27
$data[1,1] = 'ś'
$data[1,2] = 'bratem ś .'
$data[1,3] = 'świętej'
$data[1,4] = 'ś'
$data[1,5] = 'być'
$data[1,6] = 'abbreviation'
$data[2,0] = 29
$data[2,1] = 'p'
$data[2,2] = '. p .'
$data[2,3] = 'pamięci'
$data[2,4] = 'List_świętego_Piotra'
$data[2,5] = 'pan'
$data[2,6] = 'abbreviation'
$data[3,0] = 52
$data[3,1] = 'Pinińskich'
$data[3,2] = 'z Pinińskich właścicieli'
$data[3,3] = 'Piniński'
$data[3,4] = 'Pinińskich'
$data[3,5] = 'piniński'
$data[3,6] = 'surname'
$data[4,0] = 54
$data[4,1] = 'Dóbr'
$data[4,2] = 'właścicieli Dóbr Strużewo'
$data[4,3] = 'dobra'
$data[4,4] = 'Dobra'
$data[4,5] = 'dzbr'
$data[4,6] = 'capitalization'
$data[5,0] = 58
$data[5,1] = 'Puszczanki'
$data[5,2] = 'adlinencjami Puszczanki ,'
$data[5,3] = 'Puszczanka'
$data[5,4] = 'puszczanka'
$data[5,5] = 'puszczanka'
$data[5,6] = 'proper name'
$data[6,0] = 138
$data[6,1] = 'Floyrana'
$data[6,2] = 'i Floyrana .'
$data[6,3] = 'Floyran'
$data[6,4] = 'Floyrana'
$data[6,5] = 'floyrana'
$data[6,6] = 'y'
$data[7,0] = 150
$data[7,1] = 'Muczynowską'
$data[7,2] = 'z Muczynowską –'
$data[7,3] = 'Muczynowska'
$data[7,4] = 'Muczynowską'
$data[7,5] = 'muczynowska'
$data[7,6] = 'surname'
$data[8,0] = 158
$data[8,1] = 'Rzotoławskim'
$data[8,2] = 'za Rzotoławskim –'
$data[8,3] = 'Rzotoławski'
$data[8,4] = 'Rzotoławskim'
$data[8,5] = 'rzotoławski'
$data[8,6] = 'surname'
$data[9,0] = 183
$data[9,1] = 'Niewiem'
$data[9,2] = '. Niewiem o'
$data[9,3] = 'niewiedzieć'
$data[9,4] = 'Niewiem'
$data[9,5] = 'niewiem'
$data[9,6] = 'nie'
$data[10,0] = 195
$data[10,1] = 'Czermiński'
$data[10,2] = 'Karol Czermiński –'
$data[10,3] = 'Czermiński'
$data[10,4] = 'czermiński'
$data[10,5] = 'czermiński'
$data[10,6] = 'surname'
$data[11,0] = 223
$data[11,1] = 'niepamiętam'
$data[11,2] = 'jakim niepamiętam .'
$data[11,3] = 'niepamiętać'
$data[11,4] = 'niepamiętam'
$data[11,5] = 'niepamiętam'
$data[11,6] = 'nie'
$data[12,0] = 235
$data[12,1] = 'Dziedziców'
$data[12,2] = ', Dziedziców Brzeżan'
$data[12,3] = 'dziedzic'
$data[12,4] = 'Dziedzic'
$data[12,5] = 'dziedzica'
$data[12,6] = 'capitalization'
$data[13,0] = 236
$data[13,1] = 'Brzeżan'
$data[13,2] = 'Dziedziców Brzeżan i'
$data[13,3] = 'Brzeżany'
$data[13,4] = 'brzeżanin'
$data[13,5] = 'brzeżan'
$data[13,6] = 'proper name'
$data[14,0] = 345
$data[14,1] = 'dobra'
$data[14,2] = 'rozległe dobra w'
$data[14,3] = 'dobra'
$data[14,4] = 'dobry'
$data[14,5] = 'dobry'
$data[14,6] = 'ambiguous'
$data[15,0] = 392
$data[15,1] = 'administracyę'
$data[15,2] = 'objął administracyę tych'
$data[15,3] = 'administracya'
$data[15,4] = 'administracyę'
$data[15,5] = 'administracy'
$data[15,6] = 'y'
$data[16,0] = 432
$data[16,1] = 'tem'
$data[16,2] = 'ujęty tem zaproponował'
$data[16,3] = 'to'
$data[16,4] = 'tema'
$data[16,5] = 'tem'
$data[16,6] = 'e'
$data[17,0] = 518
$data[17,1] = 'Posadzie'
$data[17,2] = 'na Posadzie Chyrowskiey'
$data[17,3] = 'posada'
$data[17,4] = 'Posada'
$data[17,5] = 'posad'
$data[17,6] = 'capitalization'
$data[18,0] = 544
$data[18,1] = 'Pierściorowskim'
$data[18,2] = 'pierwszym Pierściorowskim miała'
$data[18,3] = 'Pierściorowski'
$data[18,4] = 'Pierściorowskim'
$data[18,5] = 'pierściorowski'
$data[18,6] = 'surname'
$data[19,0] = 547
$data[19,1] = 'Maryę'
$data[19,2] = 'córkę Maryę ,'
$data[19,3] = 'Marya'
$data[19,4] = 'Maryę'
$data[19,5] = 'marya'
$data[19,6] = 'y'
$data[20,0] = 549
$data[20,1] = 'Kleofasę'
$data[20,2] = ', Kleofasę .'
$data[20,3] = 'Kleofasa'
$data[20,4] = 'Kleofasę'
$data[20,5] = 'kleofasa'
$data[20,6] = 'name'
$data[21,0] = 606
$data[21,1] = 'Asińdźka'
$data[21,2] = 'wiesz Asińdźka co'
$data[21,3] = 'Asińdźka'
$data[21,4] = 'asińdźka'
$data[21,5] = 'asińdźka'
$data[21,6] = 'name'
$data[22,0] = 648
$data[22,1] = 'Kleosię'
$data[22,2] = 'moją Kleosię ”'
$data[22,3] = 'Kleosia'
$data[22,4] = 'Kleosię'
$data[22,5] = 'kleosia'
$data[22,6] = 'name'
$data[23,0] = 694
$data[23,1] = 'Pannie'
$data[23,2] = 'wstąpi Pannie się'
$data[23,3] = 'panna'
$data[23,4] = 'Panna'
$data[23,5] = 'pani'
$data[23,6] = 'capitalization'
$data[24,0] = 729
$data[24,1] = 'Czołhanach'
$data[24,2] = 'w Czołhanach roku'
$data[24,3] = 'Czołhany'
$data[24,4] = 'Czołhanach'
$data[24,5] = 'czołhanin'
$data[24,6] = 'proper name'
$data[25,0] = 739
$data[25,1] = 'osiedli'
$data[25,2] = 'nauki osiedli na'
$data[25,3] = 'osiąść'
$data[25,4] = 'osiedle'
$data[25,5] = 'osiedle'
$data[25,6] = 'ambiguous'
$data[26,0] = 756
$data[26,1] = 'Bolechowie'
$data[26,2] = 'w Bolechowie umarła'
$data[26,3] = 'Bolechów'
$data[26,4] = 'Bolechowo'
$data[26,5] = 'bolech'
$data[26,6] = 'proper name'
$data[27,0] = 784
$data[27,1] = 'Maksymowic'
$data[27,2] = 'właściciela Maksymowic w'
$data[27,3] = 'Maksymowice'
$data[27,4] = 'Maksymowic'
$data[27,5] = 'maksymowica'
$data[27,6] = 'proper name'
$data[28,0] = 796
$data[28,1] = 'Stan'
$data[28,2] = 'za Stan .'
$data[28,3] = 'Stanisław'
$data[28,4] = 'Stan'
$data[28,5] = 'stan'
$data[28,6] = 'abbreviation'
$data[29,0] = 798
$data[29,1] = 'Aug'
$data[29,2] = '. Aug .'
$data[29,3] = 'August'
$data[29,4] = 'Aug'
$data[29,5] = 'aug'
$data[29,6] = 'abbreviation'
$data[30,0] = 814
$data[30,1] = 'Bludnikach'
$data[30,2] = 'w Bludnikach –'
$data[30,3] = 'Bludniki'
$data[30,4] = 'Bludnik'
$data[30,5] = 'Bludnik'
$data[30,6] = 'proper name'
$data[31,0] = 852
$data[31,1] = 'bombardowaniu'
$data[31,2] = 'Przy bombardowaniu Lwowa'
$data[31,3] = 'bombardowanie'
$data[31,4] = 'bombardować'
$data[31,5] = 'bombardować'
$data[31,6] = 'problematic'
$data[32,0] = 863
$data[32,1] = 'obrucona'
$data[32,2] = 'perzynę obrucona –'
$data[32,3] = 'obrucić'
$data[32,4] = 'obrucona'
$data[32,5] = 'obrucon'
$data[32,6] = 'spelling'
$data[33,0] = 884
$data[33,1] = 'niespotyka'
$data[33,2] = 'się niespotyka .'
$data[33,3] = 'niespotykać'
$data[33,4] = 'niespotyka'
$data[33,5] = 'niespotyka'
$data[33,6] = 'nie'
$data[34,0] = 907
$data[34,1] = 'dobrami'
$data[34,2] = 'wielkimi dobrami ,'
$data[34,3] = 'dobra'
$data[34,4] = 'dobro'
$data[34,5] = 'dobro'
$data[34,6] = 'ambiguous'
$data[35,0] = 934
$data[35,1] = 'arye'
$data[35,2] = 'różne arye krakowiaki'
$data[35,3] = 'arya'
$data[35,4] = 'arye'
$data[35,5] = 'arye'
$data[35,6] = 'y'
$data[36,0] = 968
$data[36,1] = 'Oyca'
$data[36,2] = 'mojego Oyca –'
$data[36,3] = 'oyciec'
$data[36,4] = 'Oyca'
$data[36,5] = 'oyca'
$data[36,6] = 'y'
$data[37,0] = 976
$data[37,1] = 'Oyca'
$data[37,2] = 'familii Oyca się'
$data[37,3] = 'oyciec'
$data[37,4] = 'Oyca'
$data[37,5] = 'oyca'
$data[37,6] = 'y'
$data[38,0] = 1115
$data[38,1] = 'Świerzawski'
$data[38,2] = 'wuj Świerzawski kuzyn'
$data[38,3] = 'Świerzawski'
$data[38,4] = 'świerzawski'
$data[38,5] = 'świerzawski'
$data[38,6] = 'surname'
$data[39,0] = 1117
$data[39,1] = 'Polanowskich'
$data[39,2] = 'kuzyn Polanowskich z'
$data[39,3] = 'Polanowski'
$data[39,4] = 'polanowski'
$data[39,5] = 'polanowski'
$data[39,6] = 'surname'
$data[40,0] = 1119
$data[40,1] = 'Bełzkiego'
$data[40,2] = 'z Bełzkiego często'
$data[40,3] = 'Bełzkie'
$data[40,4] = 'Bełzkiego'
$data[40,5] = 'bełziek'
$data[40,6] = 'proper name'
$data[41,0] = 1123
$data[41,1] = 'Moszkowie'
$data[41,2] = 'w Moszkowie z'
$data[41,3] = 'Moszków'
$data[41,4] = 'Moszek'
$data[41,5] = 'moszków'
$data[41,6] = 'proper name'
$data[42,0] = 1181
$data[42,1] = 'niepomięła'
$data[42,2] = 'Komornikowej niepomięła –'
$data[42,3] = 'niepomiąć'
$data[42,4] = 'niepomięła'
$data[42,5] = 'niepomięło'
$data[42,6] = 'nie'
$data[43,0] = 1215
$data[43,1] = 'Treterówną'
$data[43,2] = 'siostrą Treterówną -'
$data[43,3] = 'Treterówna'
$data[43,4] = 'Treterówną'
$data[43,5] = 'treterówna'
$data[43,6] = 'surname'
$data[44,0] = 1233
$data[44,1] = 'Rosyi'
$data[44,2] = 'w Rosyi –'
$data[44,3] = 'Rosya'
$data[44,4] = 'Rosyi'
$data[44,5] = 'rosej'
$data[44,6] = 'y'
$data[45,0] = 1321
$data[45,1] = 'ś'
$data[45,2] = 'mi ś .'
$data[45,3] = 'świętej'
$data[45,4] = 'ś'
$data[45,5] = 'być'
$data[45,6] = 'abbreviation'
$data[46,0] = 1323
$data[46,1] = 'p'
$data[46,2] = '. p .'
$data[46,3] = 'pamięci'
$data[46,4] = 'List_świętego_Piotra'
$data[46,5] = 'pan'
$data[46,6] = 'abbreviation'
$data[47,0] = 1329
$data[47,1] = 'nieśmiał'
$data[47,2] = 'ja nieśmiał em'
$data[47,3] = 'nieśmieć'
$data[47,4] = 'MISPARSED'
$data[47,5] = 'nieśmiał'
$data[47,6] = 'nie'
$data[48,0] = 1334
$data[48,1] = 'niebył'
$data[48,2] = 'Wreszcie niebył em'
$data[48,3] = 'niebyć'
$data[48,4] = 'MISPARSED'
$data[48,5] = 'niebył'
$data[48,6] = 'nie'
$data[49,0] = 1340
$data[49,1] = 'Oycem'
$data[49,2] = 'z Oycem na'
$data[49,3] = 'oyciec'
$data[49,4] = 'Oycem'
$data[49,5] = 'oycie'
$data[49,6] = 'y'
$data[50,0] = 1418
$data[50,1] = 'kończ'
$data[50,2] = 'do kończ życia'
$data[50,3] = 'koniec'
$data[50,4] = 'kończyć'
$data[50,5] = 'kończ'
$data[50,6] = 'spelling'
$data[51,0] = 1439
$data[51,1] = 'assekuracyi'
$data[51,2] = 'w assekuracyi armat'
$data[51,3] = 'assekuracya'
$data[51,4] = 'assekuracyi'
$data[51,5] = 'assekuracyj'
$data[51,6] = 'y'
$data[52,0] = 1497
$data[52,1] = 'skrzętnem'
$data[52,2] = 'fortuny skrzętnem gospodarstwem'
$data[52,3] = 'skrzętny'
$data[52,4] = 'skrzętnem'
$data[52,5] = 'skrzętno'
$data[52,6] = 'e'
$data[53,0] = 1509
$data[53,1] = 'tem'
$data[53,2] = 'po tem jak'
$data[53,3] = 'to'
$data[53,4] = 'tema'
$data[53,5] = 'tem'
$data[53,6] = 'e'
$data[54,0] = 1517
$data[54,1] = 'set'
$data[54,2] = 'parę set sztuk'
$data[54,3] = 'sto'
$data[54,4] = 'seta'
$data[54,5] = 'set'
$data[54,6] = 'spelling'
$data[55,0] = 1527
$data[55,1] = 'Bludniki'
$data[55,2] = 'kupił Bludniki za'
$data[55,3] = 'Bludniki'
$data[55,4] = 'Bludnik'
$data[55,5] = 'bludnik'
$data[55,6] = 'proper name'
$data[56,0] = 1604
$data[56,1] = 'któremi'
$data[56,2] = 'nad któremi obeymował'
$data[56,3] = 'który'
$data[56,4] = 'któremi'
$data[56,5] = 'któr'
$data[56,6] = 'e'
$data[57,0] = 1605
$data[57,1] = 'obeymował'
$data[57,2] = 'któremi obeymował Dziedzictwo'
$data[57,3] = 'obeymować'
$data[57,4] = 'obeymował'
$data[57,5] = 'obeymował'
$data[57,6] = 'y'
$data[58,0] = 1609
$data[58,1] = 'juryzdyksye'
$data[58,2] = 'prawną juryzdyksye sądową'
$data[58,3] = 'juryzdyksya'
$data[58,4] = 'juryzdyksye'
$data[58,5] = 'juryzdyksye'
$data[58,6] = 'y'
$data[59,0] = 1614
$data[59,1] = 'pierwszey'
$data[59,2] = 'policyjną pierwszey instantacyi'
$data[59,3] = 'pierwsza'
$data[59,4] = 'pierwszey'
$data[59,5] = 'pierwsze'
$data[59,6] = 'y'
$data[60,0] = 1615
$data[60,1] = 'instantacyi'
$data[60,2] = 'pierwszey instantacyi .'
$data[60,3] = 'instantacya'
$data[60,4] = 'instantacyi'
$data[60,5] = 'instantaka'
$data[60,6] = 'y'
$data[61,0] = 1626
$data[61,1] = 'morgów'
$data[61,2] = 'posiadanych morgów -'
$data[61,3] = 'morg'
$data[61,4] = 'mórg'
$data[61,5] = 'morge'
$data[61,6] = 'problematic'
$data[62,0] = 1734
$data[62,1] = 'Ostaszewskigo'
$data[62,2] = 'Dziadka Ostaszewskigo jak'
$data[62,3] = 'Ostaszewski'
$data[62,4] = 'Ostaszewskigo'
$data[62,5] = 'ostaszewski'
$data[62,6] = 'surname'
$data[63,0] = 1741
$data[63,1] = 'kmiecie'
$data[63,2] = 'jego kmiecie na'
$data[63,3] = 'kmieć'
$data[63,4] = 'kmieci'
$data[63,5] = 'kmieta'
$data[63,6] = 'unidentified'
$data[64,0] = 1754
$data[64,1] = 'bydle'
$data[64,2] = 'padło bydle robocze'
$data[64,3] = 'bydlę'
$data[64,4] = 'bydło'
$data[64,5] = 'byska'
$data[64,6] = 'spelling'
$data[65,0] = 1772
$data[65,1] = 'Dłużanie'
$data[65,2] = 'raz Dłużanie żalili'
$data[65,3] = 'Dłużanin'
$data[65,4] = 'Dłużanie'
$data[65,5] = 'dłużać'
$data[65,6] = 'proper name'
$data[66,0] = 1777
$data[66,1] = 'niemają'
$data[66,2] = 'że niemają dogodnego'
$data[66,3] = 'niemieć'
$data[66,4] = 'niemaja'
$data[66,5] = 'niema'
$data[66,6] = 'nie'
$data[67,0] = 1790
$data[67,1] = 'suchey'
$data[67,2] = 'morgów suchey łąki'
$data[67,3] = 'suchy'
$data[67,4] = 'suchey'
$data[67,5] = 'suchey'
$data[67,6] = 'y'
$data[68,0] = 1792
$data[68,1] = 'podedworem'
$data[68,2] = 'łąki podedworem a'
$data[68,3] = 'podedwór'
$data[68,4] = 'podedworem'
$data[68,5] = 'podedwor'
$data[68,6] = 'spelling'
$data[69,0] = 1794
$data[69,1] = 'nayskładniey'
$data[69,2] = 'a nayskładniey położone'
$data[69,3] = 'składnie'
$data[69,4] = 'nayskładniey'
$data[69,5] = 'nayskładnie'
$data[69,6] = 'y'
$data[70,0] = 1804
$data[70,1] = 'Niebył'
$data[70,2] = '. Niebył to'
$data[70,3] = 'niebyć'
$data[70,4] = 'Niebyła'
$data[70,5] = 'niebył'
$data[70,6] = 'nie'
$data[71,0] = 1830
$data[71,1] = 'warżenia'
$data[71,2] = 'prawo warżenia piwa'
$data[71,3] = 'warżyć'
$data[71,4] = 'warżenia'
$data[71,5] = 'warżenie'
$data[71,6] = 'spelling'
$data[72,0] = 1881
$data[72,1] = 'miarkmi'
$data[72,2] = '– miarkmi wianki'
$data[72,3] = 'miarka'
$data[72,4] = 'miarkmi'
$data[72,5] = 'miarko'
$data[72,6] = 'spelling'
$data[73,0] = 1913
$data[73,1] = 'takiem'
$data[73,2] = 'przy takiem obciążeniu'
$data[73,3] = 'taki'
$data[73,4] = 'takiem'
$data[73,5] = 'tak'
$data[73,6] = 'e'
$data[74,0] = 1938
$data[74,1] = 'jurysdykcyi'
$data[74,2] = 'urzędowa jurysdykcyi był'
$data[74,3] = 'jurysdykcya'
$data[74,4] = 'jurysdykcyi'
$data[74,5] = 'jurysdykcyj'
$data[74,6] = 'y'
$data[75,0] = 1941
$data[75,1] = 'Dziedzica'
$data[75,2] = 'Herb Dziedzica –'
$data[75,3] = 'dziedzic'
$data[75,4] = 'Dziedzic'
$data[75,5] = 'dziedzica'
$data[75,6] = 'capitalization'
$data[76,0] = 1947
$data[76,1] = 'Bludniki'
$data[76,2] = 'Dominium Bludniki –'
$data[76,3] = 'Bludniki'
$data[76,4] = 'Bludnik'
$data[76,5] = 'bludnik'
$data[76,6] = 'proper name'
$data[77,0] = 1952
$data[77,1] = 't'
$data[77,2] = 'i t .'
$data[77,3] = 'tym'
$data[77,4] = 'tona'
$data[77,5] = 't'
$data[77,6] = 'abbreviation'
$data[78,0] = 1954
$data[78,1] = 'p'
$data[78,2] = '. p .'
$data[78,3] = 'podobne'
$data[78,4] = 'List_świętego_Piotra'
$data[78,5] = 'pan'
$data[78,6] = 'abbreviation'
$data[79,0] = 1957
$data[79,1] = 'niemógł'
$data[79,2] = 'Nieszlachcic niemógł kupować'
$data[79,3] = 'niemóc'
$data[79,4] = 'niemógł'
$data[79,5] = 'niemógł'
$data[79,6] = 'nie'
$data[80,0] = 1966
$data[80,1] = 'mieycus'
$data[80,2] = 'swojem mieycus Mandatariusza'
$data[80,3] = 'mieysce'
$data[80,4] = 'mieycus'
$data[80,5] = 'mieycus'
$data[80,6] = 'spelling'
$data[81,0] = 1997
$data[81,1] = 'niemogło'
$data[81,2] = '– niemogło się'
$data[81,3] = 'niemóc'
$data[81,4] = 'niemogło'
$data[81,5] = 'niemogło'
$data[81,6] = 'nie'
$data[82,0] = 2013
$data[82,1] = 'nieprzyniósł'
$data[82,2] = 'Dominii nieprzyniósł kartki'
$data[82,3] = 'nieprzynieść'
$data[82,4] = 'nieprzyniósł'
$data[82,5] = 'nieprzyniósło'
$data[82,6] = 'nie'
$data[83,0] = 2121
$data[83,1] = 'ludowemi'
$data[83,2] = 'dziś ludowemi zwanych'
$data[83,3] = 'ludowe'
$data[83,4] = 'ludowemi'
$data[83,5] = 'ludowa'
$data[83,6] = 'e'
$data[84,0] = 2149
$data[84,1] = 'manipulacyi'
$data[84,2] = '– manipulacyi urzędowey'
$data[84,3] = 'manipulacya'
$data[84,4] = 'manipulacyi'
$data[84,5] = 'manipulaca'
$data[84,6] = 'y'
$data[85,0] = 2150
$data[85,1] = 'urzędowey'
$data[85,2] = 'manipulacyi urzędowey –'
$data[85,3] = 'urzędowa'
$data[85,4] = 'urzędowey'
$data[85,5] = 'urzędowie'
$data[85,6] = 'y'
$data[86,0] = 2183
$data[86,1] = 'Dziedzica'
$data[86,2] = 'od Dziedzica pensyę'
$data[86,3] = 'dziedzic'
$data[86,4] = 'Dziedzic'
$data[86,5] = 'dziedzica'
$data[86,6] = 'capitalization'
$data[87,0] = 2196
$data[87,1] = 'ordynaryi'
$data[87,2] = 'korcy ordynaryi –'
$data[87,3] = 'ordynarya'
$data[87,4] = 'ordynaryi'
$data[87,5] = 'ordynary'
$data[87,6] = 'y'
$data[88,0] = 2211
$data[88,1] = 'Dziedzica'
$data[88,2] = 'zastępował Dziedzica .'
$data[88,3] = 'dziedzic'
$data[88,4] = 'Dziedzic'
$data[88,5] = 'dziedzica'
$data[88,6] = 'capitalization'
$data[89,0] = 2215
$data[89,1] = 'szczupłey'
$data[89,2] = 'tak szczupłey dotacyi'
$data[89,3] = 'szczupła'
$data[89,4] = 'szczupłey'
$data[89,5] = 'szczupłea'
$data[89,6] = 'y'
$data[90,0] = 2216
$data[90,1] = 'dotacyi'
$data[90,2] = 'szczupłey dotacyi uwzględniając'
$data[90,3] = 'dotacya'
$data[90,4] = 'dotacyi'
$data[90,5] = 'dotacy'
$data[90,6] = 'y'
$data[91,0] = 2256
$data[91,1] = 'przytem'
$data[91,2] = 'mając przytem jakiś'
$data[91,3] = 'przyto'
$data[91,4] = 'przytem'
$data[91,5] = 'przyt'
$data[91,6] = 'e'
$data[92,0] = 2274
$data[92,1] = 'lepiey'
$data[92,2] = 'nieraz lepiey jak'
$data[92,3] = 'dobrze'
$data[92,4] = 'lepiey'
$data[92,5] = 'lepiey'
$data[92,6] = 'y'
$data[93,0] = 2298
$data[93,1] = 'ładąn'
$data[93,2] = '– ładąn parą'
$data[93,3] = 'ładna'
$data[93,4] = 'ładąn'
$data[93,5] = 'ładąn'
$data[93,6] = 'spelling'
$data[94,0] = 2328
$data[94,1] = 'Bludnikach'
$data[94,2] = 'w Bludnikach –'
$data[94,3] = 'Bludniki'
$data[94,4] = 'Bludnik'
$data[94,5] = 'Bludnik'
$data[94,6] = 'proper name'
$data[95,0] = 2331
$data[95,1] = 'Bludnik'
$data[95,2] = 'do Bludnik należały'
$data[95,3] = 'Bludniki'
$data[95,4] = 'Bludnik'
$data[95,5] = 'bludnik'
$data[95,6] = 'proper name'
$data[96,0] = 2335
$data[96,1] = 'Siedliska'
$data[96,2] = 'i Siedliska –'
$data[96,3] = 'Siedliska'
$data[96,4] = 'siedliski'
$data[96,5] = 'siedliska'
$data[96,6] = 'proper name'
$data[97,0] = 2340
$data[97,1] = 'Dziedzica'
$data[97,2] = 'nowego Dziedzica –'
$data[97,3] = 'dziedzic'
$data[97,4] = 'Dziedzic'
$data[97,5] = 'dziedzica'
$data[97,6] = 'capitalization'
$data[98,0] = 2403
$data[98,1] = 'zkończone'
$data[98,2] = 'białą zkończone )'
$data[98,3] = 'zkończyć'
$data[98,4] = 'zkończone'
$data[98,5] = 'zkończone'
$data[98,6] = 'spelling'
$data[99,0] = 2415
$data[99,1] = 'niewiem'
$data[99,2] = '– niewiem czy'
$data[99,3] = 'niewiem'
$data[99,4] = 'Niewiem'
$data[99,5] = 'niew'
$data[99,6] = 'nie'
$data[100,0] = 2450
$data[100,1] = 'człowiecze'
$data[100,2] = 'Ty człowiecze na'
$data[100,3] = 'człowiek'
$data[100,4] = 'człowieczy'
$data[100,5] = 'człowiec'
$data[100,6] = 'grammar'
$data[101,0] = 2454
$data[101,1] = 'nieporadzisz'
$data[101,2] = 'mi nieporadzisz –'
$data[101,3] = 'nieporadzić'
$data[101,4] = 'niePoradzisz'
$data[101,5] = 'nieporadzisz'
$data[101,6] = 'nie'
$data[102,0] = 2508
$data[102,1] = 'mojey'
$data[102,2] = 'za mojey pamięci'
$data[102,3] = 'mój'
$data[102,4] = 'mojey'
$data[102,5] = 'mojenie'
$data[102,6] = 'y'
$data[103,0] = 2513
$data[103,1] = 'Bludnikach'
$data[103,2] = 'w Bludnikach .'
$data[103,3] = 'Bludniki'
$data[103,4] = 'Bludnik'
$data[103,5] = 'Bludnik'
$data[103,6] = 'proper name'
$data[104,0] = 2564
$data[104,1] = 'późney'
$data[104,2] = 'tak późney dokupił'
$data[104,3] = 'późno'
$data[104,4] = 'późney'
$data[104,5] = 'późnoy'
$data[104,6] = 'y'
$data[105,0] = 2573
$data[105,1] = 'Dochorowie'
$data[105,2] = 'W Dochorowie osadził'
$data[105,3] = 'Dochorów'
$data[105,4] = 'Dochorowie'
$data[105,5] = 'dochór'
$data[105,6] = 'proper name'
$data[106,0] = 2575
$data[106,1] = 'nayjstarszego'
$data[106,2] = 'osadził nayjstarszego syna'
$data[106,3] = 'stary'
$data[106,4] = 'nayjstarszego'
$data[106,5] = 'niejstarszy'
$data[106,6] = 'y'
$data[107,0] = 2589
$data[107,1] = 'Toje'
$data[107,2] = '– Toje diło'
$data[107,3] = 'toje'
$data[107,4] = 'Toje'
$data[107,5] = 'tój'
$data[107,6] = 'foreign'
$data[108,0] = 2595
$data[108,1] = 'Mychayłowu'
$data[108,2] = 'mojomu Mychayłowu .'
$data[108,3] = 'mychayłowu'
$data[108,4] = 'Mychayłowu'
$data[108,5] = 'mychaył'
$data[108,6] = 'foreign'
$data[109,0] = 2597
$data[109,1] = 'Bludniki'
$data[109,2] = '. Bludniki oddał'
$data[109,3] = 'Bludniki'
$data[109,4] = 'Bludnik'
$data[109,5] = 'bludnik'
$data[109,6] = 'proper name'
$data[110,0] = 2601
$data[110,1] = 'iOycu'
$data[110,2] = 'Józefowi iOycu mojemu'
$data[110,3] = 'ioyciec'
$data[110,4] = 'iOycu'
$data[110,5] = 'jOyc'
$data[110,6] = 'y'
$data[111,0] = 2628
$data[111,1] = 'Oyca'
$data[111,2] = 'mojego Oyca pono'
$data[111,3] = 'oyciec'
$data[111,4] = 'Oyca'
$data[111,5] = 'oyca'
$data[111,6] = 'y'
$data[112,0] = 2642
$data[112,1] = 'Siemginowa'
$data[112,2] = 'do Siemginowa –'
$data[112,3] = 'Siemginów'
$data[112,4] = 'Siemginowa'
$data[112,5] = 'siemginowa'
$data[112,6] = 'proper name'
$data[113,0] = 2646
$data[113,1] = 'dośmierci'
$data[113,2] = 'aż dośmierci Oyca'
$data[113,3] = 'dośmierć'
$data[113,4] = 'dośmierci'
$data[113,5] = 'dośmiert'
$data[113,6] = 'spelling'
$data[114,0] = 2647
$data[114,1] = 'Oyca'
$data[114,2] = 'dośmierci Oyca przy'
$data[114,3] = 'oyciec'
$data[114,4] = 'Oyca'
$data[114,5] = 'oyca'
$data[114,6] = 'y'
$data[115,0] = 2710
$data[115,1] = 'niewinem'
$data[115,2] = 'Lipa niewinem .'
$data[115,3] = 'niewiedzieć'
$data[115,4] = 'Niewino'
$data[115,5] = 'niewine'
$data[115,6] = 'nie'
$data[116,0] = 2720
$data[116,1] = 'Bludnikami'
$data[116,2] = 'z Bludnikami o'
$data[116,3] = 'Bludniki'
$data[116,4] = 'Bludnik'
$data[116,5] = 'bludnik'
$data[116,6] = 'proper name'
$data[117,0] = 2728
$data[117,1] = 'Temerowiec'
$data[117,2] = 'z Temerowiec ,'
$data[117,3] = 'Temerowice'
$data[117,4] = 'Temerowiec'
$data[117,5] = 'temerowiec'
$data[117,6] = 'proper name'
$data[118,0] = 2751
$data[118,1] = 'Oycem'
$data[118,2] = 'za Oycem doszedłszy'
$data[118,3] = 'oyciec'
$data[118,4] = 'Oycem'
$data[118,5] = 'oycie'
$data[118,6] = 'y'
$data[119,0] = 2764
$data[119,1] = 'sukcessyi'
$data[119,2] = 'prawem sukcessyi przeszedł'
$data[119,3] = 'sukcessya'
$data[119,4] = 'sukcessyi'
$data[119,5] = 'sukcessyj'
$data[119,6] = 'y'
$data[120,0] = 2804
$data[120,1] = 'niebył'
$data[120,2] = '– niebył gospodarzem'
$data[120,3] = 'niebyć'
$data[120,4] = 'Niebyła'
$data[120,5] = 'niebył'
$data[120,6] = 'nie'
$data[121,0] = 2806
$data[121,1] = 'więcey'
$data[121,2] = 'gospodarzem więcey fantasta'
$data[121,3] = 'dużo'
$data[121,4] = 'więcey'
$data[121,5] = 'więcea'
$data[121,6] = 'y'
$data[122,0] = 2847
$data[122,1] = 'gołey'
$data[122,2] = 'na gołey podłodze'
$data[122,3] = 'goła'
$data[122,4] = 'gołey'
$data[122,5] = 'goło'
$data[122,6] = 'y'
$data[123,0] = 2865
$data[123,1] = 'niemogli'
$data[123,2] = 'ci niemogli długo'
$data[123,3] = 'niemóc'
$data[123,4] = 'niemogel'
$data[123,5] = 'niemogli'
$data[123,6] = 'nie'
$data[124,0] = 2893
$data[124,1] = 'obeyściu'
$data[124,2] = 'na obeyściu niemożna'
$data[124,3] = 'obeyście'
$data[124,4] = 'obeyściu'
$data[124,5] = 'obeyść'
$data[124,6] = 'y'
$data[125,0] = 2894
$data[125,1] = 'niemożna'
$data[125,2] = 'obeyściu niemożna było'
$data[125,3] = 'niemóc'
$data[125,4] = 'niemożny'
$data[125,5] = 'niemożna'
$data[125,6] = 'nie'
$data[126,0] = 2931
$data[126,1] = 'Oyca'
$data[126,2] = 'mojego Oyca spadała'
$data[126,3] = 'oyciec'
$data[126,4] = 'Oyca'
$data[126,5] = 'oyca'
$data[126,6] = 'y'
$data[127,0] = 2938
$data[127,1] = 'Treterowej'
$data[127,2] = 'rodzonym Treterowej i'
$data[127,3] = 'Treterowa'
$data[127,4] = 'Treterowej'
$data[127,5] = 'treterowa'
$data[127,6] = 'surname'
$data[128,0] = 2940
$data[128,1] = 'Swieżaskiey'
$data[128,2] = 'i Swieżaskiey .'
$data[128,3] = 'Swieżaska'
$data[128,4] = 'Swieżaskiey'
$data[128,5] = 'swieżaskie'
$data[128,6] = 'surname'
$data[129,0] = 2948
$data[129,1] = 'Oyca'
$data[129,2] = 'śmierci Oyca z'
$data[129,3] = 'oyciec'
$data[129,4] = 'Oyca'
$data[129,5] = 'oyca'
$data[129,6] = 'y'
$data[130,0] = 2963
$data[130,1] = 'Stryiskim'
$data[130,2] = 'w Stryiskim w'
$data[130,3] = 'Stryiskie'
$data[130,4] = 'Stryiskim'
$data[130,5] = 'stryiski'
$data[130,6] = 'proper name'
$data[131,0] = 2967
$data[131,1] = 'Kruszelnicę'
$data[131,2] = 'czy Kruszelnicę –'
$data[131,3] = 'Kruszelnica'
$data[131,4] = 'Kruszelnicę'
$data[131,5] = 'kruszelnica'
$data[131,6] = 'proper name'
$data[132,0] = 2976
$data[132,1] = 'generacyi'
$data[132,2] = 'kilka generacyi na'
$data[132,3] = 'generacya'
$data[132,4] = 'generacyi'
$data[132,5] = 'generacyj'
$data[132,6] = 'y'
$data[133,0] = 2979
$data[133,1] = 'obeyściu'
$data[133,2] = 'jednem obeyściu było'
$data[133,3] = 'obeyście'
$data[133,4] = 'obeyściu'
$data[133,5] = 'obeyść'
$data[133,6] = 'y'
$data[134,0] = 2986
$data[134,1] = 'niemając'
$data[134,2] = 'Panny niemając widoków'
$data[134,3] = 'niemieć'
$data[134,4] = 'niemając'
$data[134,5] = 'niemający'
$data[134,6] = 'nie'
$data[135,0] = 3017
$data[135,1] = 'Bludnik'
$data[135,2] = 'sąsiedztwo Bludnik do'
$data[135,3] = 'Bludniki'
$data[135,4] = 'Bludnik'
$data[135,5] = 'bludnik'
$data[135,6] = 'proper name'
$data[136,0] = 3050
$data[136,1] = 'Siemiginowie'
$data[136,2] = 'na Siemiginowie już'
$data[136,3] = 'Siemiginów'
$data[136,4] = 'Siemiginowie'
$data[136,5] = 'siemigin'
$data[136,6] = 'proper name'
$data[137,0] = 3085
$data[137,1] = 'niema'
$data[137,2] = 'że niema sukcessorów'
$data[137,3] = 'niemieć'
$data[137,4] = 'niema'
$data[137,5] = 'niemy'
$data[137,6] = 'nie'
$data[138,0] = 3127
$data[138,1] = 'Kruszelnicy'
$data[138,2] = 'do Kruszelnicy o'
$data[138,3] = 'Kruszelnica'
$data[138,4] = 'Kruszelnicy'
$data[138,5] = 'kruszelnik'
$data[138,6] = 'proper name'
$data[139,0] = 3131
$data[139,1] = 'odległey'
$data[139,2] = 'mil odległey od'
$data[139,3] = 'odległa'
$data[139,4] = 'odległey'
$data[139,5] = 'odległea'
$data[139,6] = 'y'
$data[140,0] = 3133
$data[140,1] = 'Siemignowa'
$data[140,2] = 'od Siemignowa –'
$data[140,3] = 'Siemignów'
$data[140,4] = 'Siemignowa'
$data[140,5] = 'siemignowa'
$data[140,6] = 'proper name'
$data[141,0] = 3143
$data[141,1] = 'łania'
$data[141,2] = 'jak łania Panna'
$data[141,3] = 'łania'
$data[141,4] = 'łani'
$data[141,5] = 'łanie'
$data[141,6] = 'unidentified'
$data[142,0] = 3171
$data[142,1] = 'Siemiginowa'
$data[142,2] = 'sukcessorka Siemiginowa –'
$data[142,3] = 'Siemiginów'
$data[142,4] = 'Siemiginowa'
$data[142,5] = 'siemiginowa'
$data[142,6] = 'proper name'
$data[143,0] = 3185
$data[143,1] = 'Nawaryi'
$data[143,2] = 'koło Nawaryi i'
$data[143,3] = 'Nawarya'
$data[143,4] = 'Nawaryi'
$data[143,5] = 'nawary'
$data[143,6] = 'proper name'
$data[144,0] = 3207
$data[144,1] = 'ciepłey'
$data[144,2] = 'ale ciepłey wdowy'
$data[144,3] = 'ciepły'
$data[144,4] = 'ciepłey'
$data[144,5] = 'ciepło'
$data[144,6] = 'y'
$data[145,0] = 3216
$data[145,1] = 'Szołayskiego'
$data[145,2] = 'pana Szołayskiego młodzika'
$data[145,3] = 'Szołayski'
$data[145,4] = 'Szołayskiego'
$data[145,5] = 'szołayski'
$data[145,6] = 'surname'
$data[146,0] = 3230
$data[146,1] = 'pożyciu'
$data[146,2] = 'nim pożyciu –'
$data[146,3] = 'pożycie'
$data[146,4] = 'pożyć'
$data[146,5] = 'pożyć'
$data[146,6] = 'unidentified'
$data[147,0] = 3267
$data[147,1] = 'Lesniowic'
$data[147,2] = 'do Lesniowic przyległe'
$data[147,3] = 'Lesniowice'
$data[147,4] = 'Lesniowic'
$data[147,5] = 'lesniowice'
$data[147,6] = 'proper name'
$data[148,0] = 3269
$data[148,1] = 'Mosty'
$data[148,2] = 'przyległe Mosty .'
$data[148,3] = 'Mosty'
$data[148,4] = 'Most'
$data[148,5] = 'most'
$data[148,6] = 'proper name'

$ws.Range("A2:G150").Value2 = $data

# Apply the numeric-column style (bold/border/centered, same as A2) to the newly added rows
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A147:A150").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Restore view/selection as in the target workbook
$ws.Activate() | Out-Null
$ws.Range("G102").Select() | Out-Null
